# "add buy and sell product"
# - Sheet1!A6 was "SP00005" (a typo/placeholder) -> rename it to "KV000005"
#   so it matches the "KV######" naming used by the other region codes.
# - Append a brand-new region row (row 7): KV000006 / "Tây Nam Bộ".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "KV000005"

$ws.Range("A7").Value = "KV000006"
$ws.Range("B7").Value = "Tây Nam Bộ"
